$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.027.41'
$ws.Range('E2').Value = '  -1.90%  '

$ws.Range('D3').Value = '1.417.25'
$ws.Range('E3').Value = '  -1.86%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.85%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9997'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.79%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '276.18'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.40%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3685'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.58%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3111'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.90%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.83'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.66%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.040'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.37%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06514'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.03%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9982'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.02%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.494'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.80%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.63'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.89%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.183'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.86%  '

$ws.Range('D16').Value = '1.414.60'
$ws.Range('E16').Value = '  -2.63%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001019'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.84%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.05679'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.05%  '

$ws.Range('E19').Value = '  -0.70%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.15'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -7.48%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.602'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.75'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.59%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.02'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.62%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.238'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.74%  '

$ws.Range('D25').Value = '20.025.35'
$ws.Range('E25').Value = '  -1.86%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.282'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.50%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '133.26'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -6.68%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.23'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.53%  '

$ws.Range('D29').Value = '1.571.59'
$ws.Range('E29').Value = '  -2.74%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '110.22'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.891'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +11.15%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.214'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.55%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8147'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -10.71%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07775'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.51%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.460'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.33%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.898'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.05%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05829'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.11%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.064'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.86%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9994'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.82%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02054'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.14%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.48'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.71%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.103'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.71%  '

$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1875'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.60%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5306'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.12%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.32'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.97%  '

$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.535'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.49%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '117.45'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +6.21%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5187'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.07%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.769'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.64%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.035'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.12%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.000'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.86%  '
